$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 22; this shifts the existing rows 22-115
# down to 23-116 and copies formatting (incl. the date number format)
# from the row above, matching the canonical diff.
$ws.Rows.Item(22).Insert()

# Populate the newly inserted row 22 with the new record's data.
$ws.Range("A22").Value = 4
$ws.Range("B22").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C22").Value = "Los Lagos"
$ws.Range("D22").Value = 44659
$ws.Range("E22").Value = 10
$ws.Range("F22").Value = 100112052
$ws.Range("G22").Value = "Albahaca"
$ws.Range("H22").Value = "Sin especificar"
$ws.Range("I22").Value = "Primera"
$ws.Range("J22").Value = 90
$ws.Range("K22").Value = 5000
$ws.Range("L22").Value = 5000
$ws.Range("M22").Value = 5000
$ws.Range("N22").Value = "`$/docena de matas"
$ws.Range("O22").Value = "Región Metropolitana"
$ws.Range("P22").Value = 833
$ws.Range("Q22").Value = 6
$ws.Range("R22").Value = "Hortaliza"
